# Update the "想去人数" (interested-attendee count) figures in the
# "展览" and "全部类型" sheets to match the freshly scraped data.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 191
$ws.Range("F6").Value = 368
$ws.Range("F7").Value = 166
$ws.Range("F9").Value = 2338
$ws.Range("F10").Value = 124
$ws.Range("F13").Value = 1415
$ws.Range("F14").Value = 503
$ws.Range("F18").Value = 16
$ws.Range("F19").Value = 176
$ws.Range("F20").Value = 194
$ws.Range("F21").Value = 208
$ws.Range("F22").Value = 209
$ws.Range("F24").Value = 83
$ws.Range("F26").Value = 1457
$ws.Range("F27").Value = 13
$ws.Range("F28").Value = 368
$ws.Range("F29").Value = 203

# --- 全部类型 sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 191
$ws.Range("F7").Value = 368
$ws.Range("F8").Value = 166
$ws.Range("F10").Value = 2339
$ws.Range("F11").Value = 124
$ws.Range("F14").Value = 1415
$ws.Range("F15").Value = 503
$ws.Range("F19").Value = 16
$ws.Range("F20").Value = 176
$ws.Range("F21").Value = 194
$ws.Range("F22").Value = 208
$ws.Range("F23").Value = 209
$ws.Range("F25").Value = 83
$ws.Range("F27").Value = 1457
$ws.Range("F28").Value = 13
$ws.Range("F29").Value = 368
$ws.Range("F30").Value = 203
